$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: merge the "Requirement" run and the following
#    single-space run into one run "Requirement " (trailing space kept).
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:jc w:val="center"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Software </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Requirement </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Specification </w:t>
  </w:r>
</w:p>
'@
[void]$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) "Atomic Functional Requirements Shell (Volere) (cards)" paragraph:
#    split " (Volere)" into " (" + "Volere" (wrapped with spellStart /
#    spellEnd proofErr markers) and turn " (cards)" into " " + "cards)".
# ---------------------------------------------------------------------
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Atomic Functional Requirements Shell*") {
        $targetIndex = $i
        break
    }
}

$shellPara = $d.Paragraphs($targetIndex)
$shellXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="num" w:pos="1440"/>
    </w:tabs>
  </w:pPr>
  <w:r>
    <w:t>Atomic Functional Requirements Shell</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> (</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Volere</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>cards)</w:t>
  </w:r>
</w:p>
'@
[void]$shellPara.Range.InsertXML($shellXml)

# ---------------------------------------------------------------------
# 3) Remove the trailing empty paragraph at the very end of the body
#    (right before the sectPr).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
if ($lastPara.Range.Text.Trim() -eq "") {
    $prevPara = $d.Paragraphs($count - 1)
    $killRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
    $killRange.Delete()
}
